$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.934461759624336
$ws.Range("C2").Value = 3.343285591175459
$ws.Range("D2").Value = -2.437781837780359

$ws.Range("B3").Value = 1.167971958220894
$ws.Range("C3").Value = 1.173
$ws.Range("D3").Value = 0.5431565732135101

$ws.Range("B4").Value = 1.210393609907037
$ws.Range("C4").Value = 1.160939830827061
$ws.Range("D4").Value = 0.4730805793507092

$ws.Range("B5").Value = 1.675738519128104
$ws.Range("C5").Value = 2.600672032961413
$ws.Range("D5").Value = -0.9358099920757432
